# Add VPC figures: insert a new "Sensitivity" section (sensXls / sensSheet)
# into the Workflow sheet, above the "Following entries define tasks..." block.

$wb = $excel.ActiveWorkbook
$wsWorkflow = $wb.Worksheets.Item("Workflow")
$wsOutput = $wb.Worksheets.Item("Output")

# Insert three new rows at row 22 (pushes the Tasks section down to 25-29),
# inheriting the formatting from the row above (row 21), same as native Excel
# "Insert Copied Cells" / row insert behaviour.
$wsWorkflow.Rows("22:24").Insert()

# Row 22: new sub-section header "Sensitivity" (style matches the other
# sub-section header rows, e.g. row3 "simulation", row10 "population", ...)
$wsWorkflow.Range("A22:F22").Style = "Accent2"
$wsWorkflow.Rows("22").RowHeight = 49.8
$wsWorkflow.Range("A22").Value = ""
$wsWorkflow.Range("B22").Value = "Sensitivity"
$wsWorkflow.Range("C22").Value = ""

# Row 23: sensXls parameter
$wsWorkflow.Rows("23").RowHeight = 26.4
$wsWorkflow.Range("A23").Value = "sensXls"
$wsWorkflow.Range("B23").Value = "xlsfilefor sensitivity Parameter definition; if it is empty, sheet is in this xlsfile"

# Row 24: sensSheet parameter
$wsWorkflow.Rows("24").RowHeight = 26.4
$wsWorkflow.Range("A24").Value = "sensSheet"
$wsWorkflow.Range("B24").Value = "xlssheet for sensitivity Parameter definition; if empty first sheet is taken"

# Update the active selection / view to match the edited area.
$wsWorkflow.Range("C24").Select()
$excel.ActiveWindow.ScrollRow = 19

# Make "Workflow" the active (selected) sheet/tab, and "Output" no longer
# the tab-selected sheet.
$wsWorkflow.Activate()
